# Update automàtic: dades i banners [2026-02-20 09:45]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-20 09:45:32"
$ws.Range("H2").Value = "93%"
$ws.Range("I2").Value = "0.2 mm"
$ws.Range("J2").Value = "1020.7 hPa"
$ws.Range("K2").Value = "2.3 MJ/m2"
$ws.Range("M2").Value = "11.6 °C 9:29 TU"
$ws.Range("O2").Value = "2.8 °C"
